$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A (date-looking strings) to be stored as text rather than
# being auto-converted to a date serial number by Excel's type inference.
$ws.Range("A50:A51").NumberFormat = "@"

$ws.Range("A50").Value = "2024-05-15"
$ws.Range("B50").Value = "12:11:17"
$ws.Range("C50").Value = "-"
$ws.Range("D50").Value = "Cámara no detecta Power CP"
$ws.Range("E50").Value = "-"
$ws.Range("F50").Value = "-"
$ws.Range("G50").Value = "-"

$ws.Range("A51").Value = "2024-05-15"
$ws.Range("B51").Value = "12:11:20"
$ws.Range("C51").Value = "-"
$ws.Range("D51").Value = "Cámara no detecta foam derecho"
$ws.Range("E51").Value = "-"
$ws.Range("F51").Value = "-"
$ws.Range("G51").Value = "-"

# Reset style on the date cells back to Normal so no extra cell-level
# style reference (s="...") is left on them, matching the rest of the sheet.
$ws.Range("A50:A51").Style = "Normal"
